$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "59.255.49"
$ws.Range("E2").Value = "  -2.53%  "
Set-TextValue "D3" "2.584.16"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue "D5" "560.91"
$ws.Range("E5").Value = "  -1.50%  "
Set-TextValue "D6" "142.59"
$ws.Range("E6").Value = "  -3.73%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -1.87%  "
Set-TextValue "D9" "2.593.58"
$ws.Range("E9").Value = "  -3.50%  "
Set-TextValue "D10" "6.64"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("E12").Value = "  +11.14%  "
Set-TextValue "D13" "0.354"
$ws.Range("E13").Value = "  +2.76%  "
Set-TextValue "D14" "3.041.88"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D15" "59.216.30"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D16" "23.05"
$ws.Range("E16").Value = "  +5.10%  "
$ws.Range("E17").Value = "  -1.00%  "
Set-TextValue "D18" "2.577.56"
$ws.Range("E18").Value = "  -3.39%  "
Set-TextValue "D19" "4.57"
$ws.Range("E19").Value = "  +0.04%  "
Set-TextValue "D20" "336.95"
$ws.Range("E20").Value = "  -2.62%  "
Set-TextValue "D21" "10.37"
$ws.Range("E21").Value = "  -1.27%  "
Set-TextValue "D22" "6.43"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +0.36%  "
Set-TextValue "D24" "64.08"
$ws.Range("E24").Value = "  -4.36%  "
Set-TextValue "D25" "0.467"
$ws.Range("E25").Value = "  +5.44%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D31" "6.12"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "1.67"
$ws.Range("E32").Value = "  -3.18%  "
Set-TextValue "D33" "158.88"
$ws.Range("E33").Value = "  +2.27%  "
Set-TextValue "D34" "19.00"
$ws.Range("E34").Value = "  -1.90%  "
Set-TextValue "D35" "4.04"
$ws.Range("E35").Value = "  -2.17%  "
Set-TextValue "D36" "1.17"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D37" "0.876"
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D38" "0.870"
$ws.Range("E38").Value = "  -5.88%  "
Set-TextValue "D39" "37.52"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("E41").Value = "  -0.33%  "
Set-TextValue "D42" "292.97"
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("E43").Value = "  +5.60%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D47" "10.65"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D48" "0.0535"
$ws.Range("E48").Value = "  -3.14%  "
Set-TextValue "D49" "19.04"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("E50").Value = "  -1.10%  "
Set-TextValue "D51" "18.63"
$ws.Range("E51").Value = "  -1.87%  "
